$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the summary rows (2-5): the "uoa4"/TRY/25059631 row that was
# first among the data rows now becomes the last, while the others shift
# up one row - i.e. new expandable rows are appended below the previously
# first entry.
$ws.Range("A2").Value = "uoa1"
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = "USD"
$ws.Range("D2").Value = 16823445.68

$ws.Range("A3").Value = "uoa2"
$ws.Range("B3").Value = 2014
$ws.Range("C3").Value = "USD"
$ws.Range("D3").Value = 4005582.31

$ws.Range("A4").Value = "uoa3"
$ws.Range("B4").Value = 2023
$ws.Range("C4").Value = "TRY"
$ws.Range("D4").Value = 20098221.34

$ws.Range("A5").Value = "uoa4"
$ws.Range("B5").Value = 2023
$ws.Range("C5").Value = "TRY"
$ws.Range("D5").Value = 25059631

# Update the active selection to match the edited workbook's cursor position.
$ws.Range("H10").Select()
